# Apply the "I0 and IF added" change:
# Adds two new columns (I = "I0", J = "IF") to the sheet, with a header
# in row 1 (styled like the existing headers) and numeric data for rows 2-37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, center/top alignment) from the
# existing header cell H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows 2-37 (I, J values) ---
$values = @{
    2  = @(7, 7)
    3  = @(6, 6)
    4  = @(9, 9)
    5  = @(8, 8)
    6  = @(5, 6)
    7  = @(9, 9)
    8  = @(9, 9)
    9  = @(8, 8)
    10 = @(8, 8)
    11 = @(7, 7)
    12 = @(7, 7)
    13 = @(7, 7)
    14 = @(7, 7)
    15 = @(6, 6)
    16 = @(8, 8)
    17 = @(5, 5)
    18 = @(7, 7)
    19 = @(7, 7)
    20 = @(7, 8)
    21 = @(8, 8)
    22 = @(7, 7)
    23 = @(8, 8)
    24 = @(7, 7)
    25 = @(8, 8)
    26 = @(7, 8)
    27 = @(8, 8)
    28 = @(7, 7)
    29 = @(8, 8)
    30 = @(7, 7)
    31 = @(7, 7)
    32 = @(6, 6)
    33 = @(7, 7)
    34 = @(5, 5)
    35 = @(5, 5)
    36 = @(5, 5)
    37 = @(3, 3)
}

foreach ($row in $values.Keys | Sort-Object) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
}
